# Apply the Jul 15 2024 cryptos-list refresh (prices / 1h volume %, plus a few
# re-ordered rows) produced by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is Cell = New value. Values that are purely numeric-looking text
# (e.g. '562.46') are given a leading apostrophe so Excel keeps them as text
# instead of silently converting them to numbers, matching the original
# inlineStr/text formatting of the Price column.
$updates = [ordered]@{
    'D2' = '62.490.75'
    'E2' = '  +4.00%  '
    'D3' = '3.341.87'
    'E3' = '  +4.29%  '
    'E4' = '  +0.01%  '
    'D5' = '''562.46'
    'E5' = '  +4.73%  '
    'D6' = '''151.71'
    'E6' = '  +4.58%  '
    'E7' = '  -0.04%  '
    'D8' = '3.345.37'
    'E8' = '  +4.19%  '
    'D9' = '''0.537'
    'E9' = '  +1.49%  '
    'D10' = '''7.40'
    'E10' = '  +1.03%  '
    'E11' = '  +3.83%  '
    'E12' = '  +0.41%  '
    'D13' = '3.921.45'
    'E13' = '  +4.40%  '
    'E14' = '  +0.54%  '
    'E15' = '  +3.38%  '
    'E16' = '  +2.90%  '
    'D17' = '62.565.32'
    'E17' = '  +4.04%  '
    'D18' = '3.337.33'
    'E18' = '  +3.93%  '
    'E19' = '  +1.37%  '
    'D20' = '''13.78'
    'E20' = '  +4.90%  '
    'D21' = '''8.42'
    'E21' = '  +1.33%  '
    'D22' = '''385.01'
    'E22' = '  +2.45%  '
    'E23' = '  +0.13%  '
    'E24' = '  +1.52%  '
    'D25' = '''69.95'
    'E25' = '  -0.11%  '
    'E26' = '  +5.63%  '
    'D27' = '''9.07'
    'E27' = '  +3.44%  '
    'D28' = '0.0₃0952'
    'E28' = '  +6.08%  '
    'D29' = '''1.00'
    'E29' = '  -0.13%  '
    'D30' = '''6.60'
    'E30' = '  +7.62%  '
    'E31' = '  +4.13%  '
    'E32' = '  +3.67%  '
    'B33' = 'EthereumClassic'
    'C33' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D33' = '''22.88'
    'E33' = '  +2.34%  '
    'B34' = 'Fetch.AI'
    'C34' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D34' = '''1.31'
    'E34' = '  +9.04%  '
    'D35' = '''6.72'
    'E35' = '  +0.56%  '
    'E36' = '  +9.22%  '
    'D37' = '''159.25'
    'E37' = '  +1.69%  '
    'E38' = '  +11.88%  '
    'D39' = '''26.95'
    'E39' = '  +5.30%  '
    'D40' = '''0.0737'
    'E40' = '  +4.57%  '
    'D41' = '2.790.71'
    'E41' = '  -0.47%  '
    'D42' = '''0.0317'
    'E42' = '  +8.16%  '
    'B43' = 'Mantle'
    'C43' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D43' = '''0.744'
    'E43' = '  +3.59%  '
    'B44' = 'Filecoin'
    'C44' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D44' = '''4.25'
    'E44' = '  +0.42%  '
    'B45' = 'OKB'
    'C45' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D45' = '''40.47'
    'E45' = '  +1.44%  '
    'E46' = '  +4.82%  '
    'B47' = 'RenzoRestakedETH'
    'C47' = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
    'D47' = '3.387.25'
    'E47' = '  +4.35%  '
    'B48' = 'InjectiveProtocol'
    'C48' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D48' = '''22.03'
    'E48' = '  +6.90%  '
    'E49' = '  -1.74%  '
    'E50' = '  +2.28%  '
    'D51' = '''287.53'
    'E51' = '  +7.74%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
